$d = $word.ActiveDocument

# The "Assembly Guide" title lives in the primary page header, as its own
# paragraph/run ("Switch Adapted Talking Bluey Plush Toy" is the run above
# it). Scope all work to that header's Range so we can't touch anything
# else in the document.
$hdr = $d.Sections.Item(1).Headers.Item(1)

# Step 1: turn "Assembly Guide" into "maker Guide" (same run/formatting
# for now -- this mirrors the "Assembly" -> "maker" part of the commit).
$rng = $hdr.Range
$found = $rng.Find.Execute("Assembly Guide", $false, $false, $false, $false, $false, $true, 1, $false, "maker Guide", 2)

# Step 2: split " Guide" into its own run, re-applying the exact same
# character formatting the original "Assembly Guide" run carried (Roboto
# everywhere, bold, small/all caps, grey, 16pt) so the two runs end up
# looking identical but are distinct <w:r> elements, matching the diff.
$suffix = $hdr.Range
$found2 = $suffix.Find.Execute(" Guide")
$suffix.Font.Name = "Roboto"
$suffix.Font.NameFarEast = "Roboto"
$suffix.Font.NameBi = "Roboto"
$suffix.Font.NameOther = "Roboto"
$suffix.Font.Size = 16
$suffix.Font.Bold = $true
$suffix.Font.AllCaps = $true
$suffix.Font.Color = 0x646464

Write-Output "ReplacedTitle=$found SplitGuideRun=$found2 HeaderText=[$($hdr.Range.Text)]"
